$wb = $excel.ActiveWorkbook

# ---- Sheet: deals (sheet2) ----
$dealsWs = $wb.Worksheets.Item("deals")
$dealsWs.Range("A1").Value = "name"
$dealsWs.Range("B1").Value = "desc"
$dealsWs.Range("C1").Value = "prob"

$dealsWs.Range("A2").Value = "riyaz"
$dealsWs.Range("B2").Value = "This is a sample test case 1"
$dealsWs.Range("C2").Value = 100

$dealsWs.Range("A3").Value = "akram"
$dealsWs.Range("B3").Value = "This is a sample test case 2"
$dealsWs.Range("C3").Value = 110

$dealsWs.Range("A4").Value = "fahed"
$dealsWs.Range("B4").Value = "This is a sample test case 3"
$dealsWs.Range("C4").Value = 135

$dealsWs.Columns.Item(2).AutoFit()

# ---- Sheet: tasks (sheet3) ----
$tasksWs = $wb.Worksheets.Item("tasks")
$tasksWs.Range("B1").Value = "comp"
$tasksWs.Range("C1").Value = "iden"
$tasksWs.Range("A1").Value = "Title"

$tasksWs.Range("A2").Value = "Mr."
$tasksWs.Range("B2").Value = "done"
$tasksWs.Range("C2").Value = "class"

$tasksWs.Range("A3").Value = "mrs."
$tasksWs.Range("B3").Value = "progress"
$tasksWs.Range("C3").Value = "method"

$tasksWs.Range("A4").Value = "Dr."
$tasksWs.Range("B4").Value = "hold"
$tasksWs.Range("C4").Value = "function"

# ---- Selections / active sheet ----
$dealsWs.Range("C4").Select()

$tasksWs.Select()
$tasksWs.Range("C4").Select()

$contactsWs = $wb.Worksheets.Item("contacts")
$contactsWs.Range("C5").Select()

$tasksWs.Select()
